$wb = $excel.ActiveWorkbook

# --- Existing sheet: ValidLogin ---
$validLogin = $wb.Worksheets.Item(1)

# --- Add new sheet "InvalidLogin" right after "ValidLogin" ---
$invalidLogin = $wb.Worksheets.Add($null, $validLogin)
$invalidLogin.Name = "InvalidLogin"

# Populate the new "InvalidLogin" sheet with header + bad credentials + expected message
$invalidLogin.Range("A1").Value = "UserName"
$invalidLogin.Range("B1").Value = "Password"
$invalidLogin.Range("A2").Value = "abcd"
$invalidLogin.Range("B2").Value = "xyz"
$invalidLogin.Range("A3").Value = "admin"
$invalidLogin.Range("B3").Value = "damager"

# Zoom both sheets to 220% (was 175% on ValidLogin) and restore each sheet's selection
$validLogin.Activate()
[void]$validLogin.Range("A2").Select()
$excel.ActiveWindow.Zoom = 220

$invalidLogin.Activate()
[void]$invalidLogin.Range("B3").Select()
$excel.ActiveWindow.Zoom = 220
